# Saldo_guide.xlsx update
# - Roll every row's date in column G from 2024-04-01 (45383) to 2024-04-03 (45385)
# - Correct the D/E/H balances for a handful of rows whose reconciled totals changed
# - Leave the active selection on K9 (matches the re-saved sheetView)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump every date in the G column (rows 2-310) to the new date.
$ws.Range("G2:G310").Value = 45385

# Row-specific corrections to D (debit), E (credit) and H (balance) columns.
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("H5").Value = 0

$ws.Range("D23").Value = 1436.03
$ws.Range("E23").Value = 0
$ws.Range("H23").Value = 1436.03

$ws.Range("D27").Value = 30340.69
$ws.Range("E27").Value = 0
$ws.Range("H27").Value = 30340.69

$ws.Range("D47").Value = 9530.27
$ws.Range("H47").Value = 9530.27

$ws.Range("D61").Value = 24.48
$ws.Range("E61").Value = 0
$ws.Range("H61").Value = 24.48

$ws.Range("D64").Value = 1825.89
$ws.Range("H64").Value = 1825.89

$ws.Range("D69").Value = 0
$ws.Range("H69").Value = 0

$ws.Range("D76").Value = 0
$ws.Range("E76").Value = 0
$ws.Range("H76").Value = 0

$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 0
$ws.Range("H78").Value = 0

$ws.Range("D129").Value = 66903.149999999994
$ws.Range("H129").Value = 66903.149999999994

$ws.Range("D130").Value = 66903.13
$ws.Range("H130").Value = 66903.13

$ws.Range("D151").Value = 12099.52
$ws.Range("E151").Value = 0
$ws.Range("H151").Value = 12099.52

$ws.Range("D169").Value = 40.14
$ws.Range("H169").Value = 40.14

$ws.Range("D187").Value = 233.72
$ws.Range("H187").Value = 233.72

$ws.Range("D200").Value = 1154.46
$ws.Range("E200").Value = 0
$ws.Range("H200").Value = 1154.46

$ws.Range("D267").Value = 24.66
$ws.Range("E267").Value = 0

$ws.Range("D274").Value = 0
$ws.Range("E274").Value = 0
$ws.Range("H274").Value = 0

# Match the saved selection state (cell K9 active on the sheet).
$ws.Range("K9").Select()
